$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 (headers) - re-label columns O..Z
# ---------------------------------------------------------------------------
$ws.Range("O1").Value = "年"
$ws.Range("P1").Value = "アイテムURL"
$ws.Range("Q1").Value = "viewingDirection"
$ws.Range("R1").Value = "ID"
$ws.Range("S1").Value = "ソート用項目"
$ws.Range("T1").Value = "機械可読ドキュメント"
$ws.Range("U1").Value = "ウェブサイトURL"
$ws.Range("V1").Value = "IIIFマニフェストURI"
$ws.Range("W1").Value = "帰属"
$ws.Range("X1").Value = "コレクション"
$ws.Range("Y1").Value = "サムネイル"
$ws.Range("Z1").Value = "利用条件"

# ---------------------------------------------------------------------------
# Row 2 (values) - capture the old text up front so overwrite order doesn't
# matter, then write the new layout.
# ---------------------------------------------------------------------------
$oldP2 = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/26/full/200,151/0/default.jpg"
$oldS2 = "fbd0479b-dbb4-4eaa-95b8-f27e1c423e4b"
$oldT2 = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/hyakki/document/fbd0479b-dbb4-4eaa-95b8-f27e1c423e4b"
$oldV2 = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/8"
$oldW2 = "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse"
$oldX2 = "東京大学総合図書館 General Library in the University of Tokyo, JAPAN"
$oldY2 = "http://iiif.io/api/presentation/2#rightToLeftDirection"
$oldA2 = $ws.Range("A2").Value2

# R2/S2 collapse into a single R2 cell (ID of the item)
$ws.Range("S2").ClearContents()
$ws.Range("R2").Value = $oldS2

$ws.Range("P2").Value = $oldT2
$ws.Range("T2").Value = $oldV2
$ws.Range("W2").Value = $oldX2
$ws.Range("X2").Value = $oldA2
$ws.Range("Y2").Value = $oldP2
$ws.Range("Z2").Value = $oldW2
$ws.Range("Q2").Value = $oldY2

# ---------------------------------------------------------------------------
# Hyperlinks - drop the old set and rebuild it against the new cells.
# (Deleting while iterating a live COM collection skips entries, so match +
# delete one at a time, re-fetching the collection on every pass.)
# ---------------------------------------------------------------------------
$oldAddrs = @('$P$2', '$Q$2', '$T$2', '$U$2', '$V$2', '$W$2', '$Y$2')
foreach ($addr in $oldAddrs) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            break
        }
    }
}

$ws.Hyperlinks.Add($ws.Range("P2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/26/full/200,151/0/default.jpg", "")
$ws.Hyperlinks.Add($ws.Range("Q2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/8/manifest", "rightToLeftDirection")
$ws.Hyperlinks.Add($ws.Range("T2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/hyakki/document/fbd0479b-dbb4-4eaa-95b8-f27e1c423e4b", "")
$ws.Hyperlinks.Add($ws.Range("U2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/hyakki/", "")
$ws.Hyperlinks.Add($ws.Range("V2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/8", "")
$ws.Hyperlinks.Add($ws.Range("Y2"), "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse", "")
$ws.Hyperlinks.Add($ws.Range("Z2"), "http://iiif.io/api/presentation/2", "")

# ---------------------------------------------------------------------------
# Styles - keep the "hyperlink" look (underline + blue font) on cells that
# still carry a link, strip it from W2 (no longer a link). `Hyperlinks.Add`
# stamps its own blue/theme style on P2/Q2/T2/U2/V2/Y2/Z2, so restore the
# sheet's original hyperlink look (font 1 / style 1) by copying formats from
# a cell that already carries it, and clear W2 back to the plain look.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("W2").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Copy() | Out-Null
$ws.Range("AB1").PasteSpecial(-4122) | Out-Null
$ws.Range("AB1").Value = "refstyle"
$ws.Range("AB1").Font.Underline = 2
$ws.Range("AB1").Font.Color = 16711680

$ws.Range("AB1").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2").PasteSpecial(-4122) | Out-Null
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$ws.Range("U2").PasteSpecial(-4122) | Out-Null
$ws.Range("V2").PasteSpecial(-4122) | Out-Null
$ws.Range("Y2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").PasteSpecial(-4122) | Out-Null

$ws.Range("AB1").Clear() | Out-Null
$excel.CutCopyMode = 0
